$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.546.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4275"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3666"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07269"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.812.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.394"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.514"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06941"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008901"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("E19").Value = "  -0.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.562.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.156"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.59%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.080.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.993"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "

# Row 29
$ws.Range("E29").Value = "  -4.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.820"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08850"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.987"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7466"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.538"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.19%  "

# Row 36
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("E37").Value = "  -1.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.797"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "

# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1661"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5072"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.499"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06482"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4672"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.610"
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
